$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6943265199661255
$ws.Range("B1").Value = 1.33771276473999
$ws.Range("C1").Value = 3.883811712265015
$ws.Range("D1").Value = 2.723388433456421
$ws.Range("E1").Value = 0.5670968294143677
